# Adds two new survey rows (a "calculate" field harvest_days and a "note"
# field harvest_days_note) right after the harvest_date question (row 201)
# in the "survey" sheet of the ODK production-event workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert two blank rows at 202 (pushes old rows 202..234 down to 204..236).
$ws.Rows.Item(202).Insert()
$ws.Rows.Item(202).Insert()

# Row 202: type=calculate, name=harvest_days, calculation=...
# (columns filled in this exact order so new shared-string entries come out
# in the same sequence as the target workbook)
$ws.Range("A202").Value = "calculate"
$ws.Range("B202").Value = "harvest_days"
$ws.Range("A203").Value = "note"
$ws.Range("B203").Value = "harvest_days_note"
$ws.Range("J202").Value = 'decimal-date-time(${harvest_date}) - decimal-date-time(${planting_date})'
$ws.Range("C203").Value = 'Días desde siembra ${harvest_days}'

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the new
# extent of the survey table (now A1:L236 instead of A1:L234).
$n = $wb.Names.Item("survey!_FilterDatabase")
$n.RefersTo = "=survey!`$A`$1:`$L`$236"

# Reflect the author's final cursor/selection position on the sheet.
$ws.Activate()
$ws.Range("C204").Select()
